# Regen sval data to filter save games.
# Update columns B:G (TB, d2S, K, IP, Win, sum) for rows 2-14 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 3.900430680208489,  0.496779210170732, 0, 9.295990156953671)
    3  = @(3.230985683306322, 1.667794583268128, 26.21740644021617,  0.496779210170732, 1, 31.61296591696135)
    4  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    6  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    7  = @(3.230985683306322, 1.667794583268128, 26.21740644021617,  0.496779210170732, 1, 31.61296591696135)
    8  = @(0.01514828764759746, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 1.630207530253468)
    9  = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 6.740334628841572)
    10 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    12 = @(0.127881588408715, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 3.097945018431574)
    13 = @(0.6753301551942219, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 2.290389397800092)
    14 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
